# Update "paises" (countries) COVID data sheet + the "last updated" timestamp.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# --- Timestamp in A1 -------------------------------------------------
$ws.Range("A1").Value = "Datos actualizados a 15 de Octubre de 2020 a las 16:43"

# --- Helper: write a full data row (Country, Casos totales, Nuevos casos,
#     Casos activos, Recuperados, Muertes hoy, Casos criticos, Muertes) ----
function Set-CountryRow {
    param($row, $country, $b, $c, $d, $e, $f, $g, $h)
    $ws.Cells.Item($row, 1).Value = $country
    $ws.Cells.Item($row, 2).Value = $b
    $ws.Cells.Item($row, 3).Value = $c
    $ws.Cells.Item($row, 4).Value = $d
    $ws.Cells.Item($row, 5).Value = $e
    $ws.Cells.Item($row, 6).Value = $f
    $ws.Cells.Item($row, 7).Value = $g
    $ws.Cells.Item($row, 8).Value = $h
}

# Plain data refreshes (no re-sort / swap of neighbouring countries)
Set-CountryRow 4   "Estados Unidos"        8159163 9120 5280891 2656356 0 73 221916
Set-CountryRow 5   "India"                 7311088 6018 6384885  814855 0 37 111348
Set-CountryRow 23  "Alemania"               344089 2347  281900   52404 0 14   9785
Set-CountryRow 60  "Uzbekistan"              62278  328   59291    2469 0  4    518
Set-CountryRow 71  "Estado de Palestina"     46100  442   39585    6114 0  8    401
Set-CountryRow 111 "Uganda"                  10117   48    6725    3296 0  1     96
Set-CountryRow 122 "Cuba"                     6062   27    5678     261 0  0    123

# Portugal overtakes Costa Rica; Japon keeps its spot but gets new figures
Set-CountryRow 49  "Portugal"                93294 2101   55081   36085 0 11   2128
Set-CountryRow 50  "Costa Rica"              91780    0   54155   36491 0  0   1134
Set-CountryRow 51  "Japon"                   90710  570   83837    5227 0  8   1646

# Kenia overtakes Azerbaiyan
Set-CountryRow 73  "Kenia"                   43143  602   31508   10830 0  8    805
Set-CountryRow 74  "Azerbaiyan"              42750    0   39570    2564 0  0    616

# Birmania overtakes Bosnia y Herzegovina
Set-CountryRow 80  "Birmania"                32351 1026   14706   16880 0 33    765
Set-CountryRow 81  "Bosnia y Herzegovina"    32224  569   24356    6896 0 14    972

# Noruega overtakes Albania
Set-CountryRow 95  "Noruega"                 16050   97   11863    3909 0  1    278
Set-CountryRow 96  "Albania"                 15955    0    9762    5759 0  0    434
